$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2-9) for columns D, L, M, N, O, P, R, S
# Derived from the diff: the data rows were re-shuffled (dates/quality/volume/prices/origin)
# while identifying (A,B,C,E-K,Q,T) columns stayed constant.

$data = @{
    2 = @{ D = 44162; L = "Tercera"; M = 500; N = 15000; O = 16000; P = 15500; R = "Región de O'Higgins"; S = 1033 }
    3 = @{ D = 44176; L = "Segunda"; M = 500; N = 15000; O = 16000; P = 15500; R = "Región Metropolitana"; S = 1033 }
    4 = @{ D = 44194; L = "Segunda"; M = 300; N = 15000; O = 16000; P = 15500; R = "Región Metropolitana"; S = 1033 }
    5 = @{ D = 44169; L = "Segunda"; M = 500; N = 15000; O = 16000; P = 15500; R = "Región de O'Higgins"; S = 1033 }
    6 = @{ D = 44159; L = "Tercera"; M = 400; N = 15500; O = 16000; P = 15750; R = "Región de O'Higgins"; S = 1050 }
    7 = @{ D = 44166; L = "Segunda"; M = 600; N = 16000; O = 17000; P = 16500; R = "Región de O'Higgins"; S = 1100 }
    8 = @{ D = 44187; L = "Primera"; M = 350; N = 16000; O = 16000; P = 16000; R = "Región Metropolitana"; S = 1067 }
    9 = @{ D = 44187; L = "Segunda"; M = 300; N = 13000; O = 13000; P = 13000; R = "Región Metropolitana"; S = 867 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
